# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the job sheets with the latest fetched values.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1883.1666
$ws.Range("I96").Value = 933
$ws.Range("J96").Value = 2833.3333
$ws.Range("K96").Value = 2799
$ws.Range("L96").Value = 8499.999899999999
$ws.Range("M96").Value = -1426
$ws.Range("N96").Value = -11245.9999
$ws.Range("H107").Value = 2796.4285
$ws.Range("I107").Value = 2796.4285
$ws.Range("K107").Value = 2796.4285
$ws.Range("M107").Value = -876.4285
$ws.Range("H135").Value = 651.6842
$ws.Range("I135").Value = 571.8
$ws.Range("K135").Value = 5146.2
$ws.Range("M135").Value = -2611.2
$ws.Range("H138").Value = 2650.4656
$ws.Range("I138").Value = 2438.75
$ws.Range("K138").Value = 7316.25
$ws.Range("M138").Value = -2176.25
$ws.Range("H141").Value = 982.6774
$ws.Range("I141").Value = 982.6774
$ws.Range("K141").Value = 2948.0322
$ws.Range("M141").Value = 2231.9678

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 316.33334
$ws.Range("I4").Value = 383.6
$ws.Range("K4").Value = 383.6
$ws.Range("M4").Value = -267.6
$ws.Range("H5").Value = 151.18518
$ws.Range("I5").Value = 263.8
$ws.Range("J5").Value = 125.59091
$ws.Range("K5").Value = 263.8
$ws.Range("L5").Value = 125.59091
$ws.Range("M5").Value = -151.8
$ws.Range("N5").Value = -349.59091
$ws.Range("H97").Value = 7464.8335
$ws.Range("I97").Value = 8197.929
$ws.Range("K97").Value = 8197.929
$ws.Range("M97").Value = -7701.929

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 151.18518
$ws.Range("I4").Value = 263.8
$ws.Range("J4").Value = 125.59091
$ws.Range("K4").Value = 263.8
$ws.Range("L4").Value = 125.59091
$ws.Range("M4").Value = -148.8
$ws.Range("N4").Value = -355.59091
$ws.Range("H105").Value = 137356
$ws.Range("I105").Value = 4048.2727
$ws.Range("J105").Value = 503952.25
$ws.Range("K105").Value = 4048.2727
$ws.Range("L105").Value = 503952.25
$ws.Range("M105").Value = -2301.2727
$ws.Range("N105").Value = -507446.25
$ws.Range("H107").Value = 1641.0358
$ws.Range("I107").Value = 980
$ws.Range("K107").Value = 980
$ws.Range("M107").Value = 940
$ws.Range("H134").Value = 3035.8538
$ws.Range("I134").Value = 2855.0286
$ws.Range("K134").Value = 8565.085800000001
$ws.Range("M134").Value = -6030.085800000001

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 125.23529
$ws.Range("I7").Value = 90.44444
$ws.Range("K7").Value = 90.44444
$ws.Range("M7").Value = 22.55556
$ws.Range("H31").Value = 3141.5715
$ws.Range("I31").Value = 2042
$ws.Range("J31").Value = 4241.143
$ws.Range("K31").Value = 2042
$ws.Range("L31").Value = 4241.143
$ws.Range("M31").Value = -1747
$ws.Range("N31").Value = -4831.143
$ws.Range("H34").Value = 3141.5715
$ws.Range("I34").Value = 2042
$ws.Range("J34").Value = 4241.143
$ws.Range("K34").Value = 2042
$ws.Range("L34").Value = 4241.143
$ws.Range("M34").Value = -1840
$ws.Range("N34").Value = -4645.143
$ws.Range("H41").Value = 9739.6
$ws.Range("I41").Value = 4674.75
$ws.Range("K41").Value = 4674.75
$ws.Range("M41").Value = -4246.75
$ws.Range("H47").Value = 39999
$ws.Range("J47").Value = 39999
$ws.Range("L47").Value = 39999
$ws.Range("N47").Value = -41131
$ws.Range("H86").Value = 5997.9
$ws.Range("I86").Value = 4997
$ws.Range("K86").Value = 4997
$ws.Range("M86").Value = -3874
$ws.Range("H89").Value = 5997.9
$ws.Range("I89").Value = 4997
$ws.Range("K89").Value = 24985
$ws.Range("M89").Value = -19369
$ws.Range("H94").Value = 1399.3334
$ws.Range("I94").Value = 1399
$ws.Range("K94").Value = 1399
$ws.Range("M94").Value = -948
$ws.Range("H132").Value = 1909.5385
$ws.Range("I132").Value = 1719.1818
$ws.Range("K132").Value = 5157.5454
$ws.Range("M132").Value = -2627.5454

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3963
$ws.Range("I62").Value = 3963
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 11889
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -11203
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3963
$ws.Range("I65").Value = 3963
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 35667
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -32235
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 480.82608
$ws.Range("I113").Value = 497.84616
$ws.Range("J113").Value = 458.7
$ws.Range("K113").Value = 1493.53848
$ws.Range("L113").Value = 1376.1
$ws.Range("M113").Value = 676.4615200000001
$ws.Range("N113").Value = -5716.1
$ws.Range("H117").Value = 2960.6667
$ws.Range("J117").Value = 6877
$ws.Range("L117").Value = 20631
$ws.Range("N117").Value = -27515

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4655
$ws.Range("I70").Value = 4299.1665
$ws.Range("J70").Value = 5366.6665
$ws.Range("K70").Value = 4299.1665
$ws.Range("L70").Value = 5366.6665
$ws.Range("M70").Value = -4029.1665
$ws.Range("N70").Value = -5906.6665
$ws.Range("H73").Value = 4655
$ws.Range("I73").Value = 4299.1665
$ws.Range("J73").Value = 5366.6665
$ws.Range("K73").Value = 4299.1665
$ws.Range("L73").Value = 5366.6665
$ws.Range("M73").Value = -3363.1665
$ws.Range("N73").Value = -7238.6665
$ws.Range("H126").Value = 5799.5625
$ws.Range("I126").Value = 5367.6665
$ws.Range("K126").Value = 16102.9995
$ws.Range("M126").Value = -13632.9995

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 30000
$ws.Range("J21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("N21").Value = -30348
$ws.Range("H22").Value = 63146.11
$ws.Range("I22").Value = 124204.22
$ws.Range("J22").Value = 2088
$ws.Range("K22").Value = 124204.22
$ws.Range("L22").Value = 2088
$ws.Range("M22").Value = -123909.22
$ws.Range("N22").Value = -2678
$ws.Range("H27").Value = 63146.11
$ws.Range("I27").Value = 124204.22
$ws.Range("J27").Value = 2088
$ws.Range("K27").Value = 124204.22
$ws.Range("L27").Value = 2088
$ws.Range("M27").Value = -124097.22
$ws.Range("N27").Value = -2302
$ws.Range("H38").Value = 100000
$ws.Range("J38").Value = 100000
$ws.Range("L38").Value = 100000
$ws.Range("N38").Value = -100820
$ws.Range("H132").Value = 46982.816
$ws.Range("I132").Value = 48501.42
$ws.Range("K132").Value = 145504.26
$ws.Range("M132").Value = -142974.26

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 31664.854
$ws.Range("I126").Value = 43649.207
$ws.Range("J126").Value = 2902.4
$ws.Range("K126").Value = 130947.621
$ws.Range("L126").Value = 8707.200000000001
$ws.Range("M126").Value = -128477.621
$ws.Range("N126").Value = -13647.2
$ws.Range("H132").Value = 59719.945
$ws.Range("I132").Value = 59719.945
$ws.Range("K132").Value = 179159.835
$ws.Range("M132").Value = -176629.835
